$d = $word.ActiveDocument

# --- Change 1 ------------------------------------------------------------
# "Lekcję 3 należy zrealizować innego dnia, niż lekcje 1 i 2."
#   -> "Lekcję 2 warto jest zrealizować tego samego dnia co lekcję 1."
$d.Content.Find.Execute(
    "Lekcję 3 należy zrealizować innego dnia, niż lekcje 1 i 2.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Lekcję 2 warto jest zrealizować tego samego dnia co lekcję 1.", 2)

# --- Change 2 ------------------------------------------------------------
# The paragraph "Lekcję 3, ze względu na jej ewaluacyjny charakter, należy
# potraktować podobnie do sprawdzianu." becomes two paragraphs:
#   "Lekcję 4 należy zrealizować innego dnia, niż lekcje 1, 2 i 3."
#   "Lekcję 4, ze względu na jej ewaluacyjny charakter, należy potraktować
#    podobnie do sprawdzianu."
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*ewaluacyjny charakter*") {
        $target = $para
        break
    }
}
$idx = $target.Index

# Renumber the lesson reference in this bullet from 3 to 4 (scoped to this
# paragraph only, so other "3"s in the document stay untouched).
$target.Range.Find.Execute("3", $true, $false, $false, $false, $false, $true, 1, $false, "4", 2)

# Find the split point: right before the ", ze względu" that keeps going with
# the original sentence.
$ip = $d.Paragraphs($idx).Range
$ip.Find.Execute(", ze względu")
$ip.Collapse(1)
$ip.InsertParagraphAfter()

# The first (new) paragraph now ends right after "Lekcję 4" -- append the new
# explanation about needing a different day than lessons 1, 2 and 3.
$firstPara = $d.Paragraphs($idx)
$firstPara.Range.InsertAfter(" należy zrealizować innego dnia, niż lekcje 1, 2 i 3.")

# The second paragraph still starts with ", ze względu ..." -- prefix it with
# "Lekcję 4" to restore a full sentence.
$secondPara = $d.Paragraphs($idx + 1)
$secondPara.Range.InsertBefore("Lekcję 4")
